$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Плечи / 0-3): C changes 0 -> 1
$ws.Range("C2").Value = 1

# Row 3 (Таз / 0-4): C changes 1 -> 3
$ws.Range("C3").Value = 3

# Row 4 (Локти / 0-2): C changes 3 -> 2, D changes 5 -> 8
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 8

# Row 5 (Колени / 0-2): C changes 2 -> 8, D changes 2 -> 4
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 4

# Row 6 (Шея / 0-8): C changes 2 -> 12, D changes 4 -> 7
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 7

# Row 7 (Осанка / 0-8): C changes 2 -> 4, D changes 13 -> 7
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 7

# Row 8 (Голень/Стопа): label "0-9" -> "100 - 110"; D changes 112 -> 109
$ws.Range("B8").Value = "100 - 110"
$ws.Range("D8").Value = 109

# Row 9 (Руки вперед): label "0-6" -> "80 - 90"; C changes 71 -> 82, D changes 65 -> 79
$ws.Range("B9").Value = "80 - 90"
$ws.Range("C9").Value = 82
$ws.Range("D9").Value = 79

# Row 10 (Руки в стороны): label "0-6" -> "85 - 90"; C changes 84 -> 81, D changes 90 -> 94
$ws.Range("B10").Value = "85 - 90"
$ws.Range("C10").Value = 81
$ws.Range("D10").Value = 94

# Row 11 (Руки подняты): label "0-5" -> "170 - 180"; C changes 217 -> 169, D changes 166 -> 187
$ws.Range("B11").Value = "170 - 180"
$ws.Range("C11").Value = 169
$ws.Range("D11").Value = 187
